{"js": "/* Word JS API (Office.js) edit script.\n * Updates the worksheet title date and every multiplication problem in the\n * 20x5 practice table. Cell/paragraph text is replaced in place via\n * Range.insertText(..., Word.InsertLocation.replace) so existing run/paragraph\n * formatting (fonts, size, alignment) is preserved; only the literal text\n * content changes, matching the source diff exactly. Because a couple of the\n * original problems repeat verbatim (e.g. \"97\u00d738=\" appears twice but maps to\n * two different results), the update is positional (row-major over the\n * table), not a global text search/replace.\n */\n\nconst newDate = \"2023-04-04 Tuesday\";\nconst newGrid = [\n  [\"67\u00d729=\", \"52\u00d723=\", \"67\u00d780=\", \"17\u00d736=\", \"70\u00d759=\"],\n  [\"46\u00d747=\", \"93\u00d777=\", \"34\u00d776=\", \"16\u00d761=\", \"64\u00d741=\"],\n  [\"48\u00d777=\", \"58\u00d749=\", \"76\u00d793=\", \"23\u00d792=\", \"78\u00d747=\"],\n  [\"88\u00d722=\", \"51\u00d733=\", \"76\u00d754=\", \"56\u00d742=\", \"83\u00d783=\"],\n  [\"65\u00d738=\", \"20\u00d779=\", \"69\u00d756=\", \"12\u00d736=\", \"83\u00d729=\"],\n  [\"94\u00d725=\", \"11\u00d757=\", \"26\u00d734=\", \"99\u00d734=\", \"50\u00d750=\"],\n  [\"72\u00d735=\", \"15\u00d730=\", \"52\u00d727=\", \"47\u00d779=\", \"45\u00d751=\"],\n  [\"51\u00d784=\", \"19\u00d785=\", \"48\u00d712=\", \"13\u00d712=\", \"14\u00d766=\"],\n  [\"38\u00d772=\", \"86\u00d728=\", \"33\u00d799=\", \"95\u00d782=\", \"18\u00d734=\"],\n  [\"41\u00d760=\", \"53\u00d714=\", \"89\u00d710=\", \"48\u00d730=\", \"71\u00d796=\"],\n  [\"92\u00d731=\", \"77\u00d764=\", \"18\u00d773=\", \"31\u00d783=\", \"85\u00d747=\"],\n  [\"61\u00d737=\", \"44\u00d779=\", \"73\u00d755=\", \"59\u00d769=\", \"88\u00d767=\"],\n  [\"57\u00d787=\", \"27\u00d731=\", \"13\u00d799=\", \"68\u00d793=\", \"85\u00d766=\"],\n  [\"26\u00d730=\", \"46\u00d789=\", \"99\u00d757=\", \"62\u00d748=\", \"90\u00d782=\"],\n  [\"54\u00d791=\", \"77\u00d762=\", \"100\u00d761=\", \"50\u00d784=\", \"49\u00d742=\"],\n  [\"49\u00d714=\", \"93\u00d765=\", \"14\u00d718=\", \"22\u00d752=\", \"53\u00d757=\"],\n  [\"20\u00d748=\", \"92\u00d779=\", \"74\u00d722=\", \"44\u00d713=\", \"14\u00d759=\"],\n  [\"80\u00d790=\", \"73\u00d754=\", \"89\u00d786=\", \"20\u00d745=\", \"40\u00d736=\"],\n  [\"65\u00d777=\", \"45\u00d764=\", \"34\u00d756=\", \"33\u00d799=\", \"83\u00d759=\"],\n  [\"11\u00d744=\", \"35\u00d748=\", \"64\u00d732=\", \"93\u00d793=\", \"52\u00d723=\"],\n];\n\n// --- Update the title paragraph (the date line above the table). ---\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length > 0) {\n  const titleRange = paragraphs.items[0].getRange();\n  titleRange.insertText(newDate, Word.InsertLocation.replace);\n}\n\n// --- Update every cell of the multiplication table, row by row. ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount, values\");\nawait context.sync();\n\nconst columnCount = table.values.length > 0 ? table.values[0].length : 0;\nif (table.rowCount !== newGrid.length || columnCount !== newGrid[0].length) {\n  throw new Error(\n    `Unexpected table shape: got ${table.rowCount}x${columnCount}, ` +\n      `expected ${newGrid.length}x${newGrid[0].length}`\n  );\n}\n\nfor (let r = 0; r < newGrid.length; r++) {\n  const row = newGrid[r];\n  for (let c = 0; c < row.length; c++) {\n    const cell = table.getCellOrNullObject(r, c);\n    const cellRange = cell.body.paragraphs.getFirstOrNullObject().getRange();\n    cellRange.insertText(row[c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Updates the title date paragraph and every multiplication problem in the\n# 20x5 practice table. Cell/paragraph Range.Text is assigned in place so the\n# existing run/paragraph formatting (fonts, size, alignment) is preserved and\n# only the literal text changes -- matching the source diff exactly. The\n# update is positional (row major over the table) rather than a global\n# text search/replace because a couple of the original problems repeat\n# verbatim (e.g. \"97\u00d738=\" appears twice but maps to two different results).\n\n$d = $word.ActiveDocument\n\n# --- Update the title paragraph (the date line above the table). ---\n$d.Paragraphs(1).Range.Text = \"2023-04-04 Tuesday\"\n\n# --- Update every cell of the multiplication table, row by row. ---\n$t = $d.Tables(1)\n\n$newValues = @(\n  ,@(\"67\u00d729=\", \"52\u00d723=\", \"67\u00d780=\", \"17\u00d736=\", \"70\u00d759=\")\n  ,@(\"46\u00d747=\", \"93\u00d777=\", \"34\u00d776=\", \"16\u00d761=\", \"64\u00d741=\")\n  ,@(\"48\u00d777=\", \"58\u00d749=\", \"76\u00d793=\", \"23\u00d792=\", \"78\u00d747=\")\n  ,@(\"88\u00d722=\", \"51\u00d733=\", \"76\u00d754=\", \"56\u00d742=\", \"83\u00d783=\")\n  ,@(\"65\u00d738=\", \"20\u00d779=\", \"69\u00d756=\", \"12\u00d736=\", \"83\u00d729=\")\n  ,@(\"94\u00d725=\", \"11\u00d757=\", \"26\u00d734=\", \"99\u00d734=\", \"50\u00d750=\")\n  ,@(\"72\u00d735=\", \"15\u00d730=\", \"52\u00d727=\", \"47\u00d779=\", \"45\u00d751=\")\n  ,@(\"51\u00d784=\", \"19\u00d785=\", \"48\u00d712=\", \"13\u00d712=\", \"14\u00d766=\")\n  ,@(\"38\u00d772=\", \"86\u00d728=\", \"33\u00d799=\", \"95\u00d782=\", \"18\u00d734=\")\n  ,@(\"41\u00d760=\", \"53\u00d714=\", \"89\u00d710=\", \"48\u00d730=\", \"71\u00d796=\")\n  ,@(\"92\u00d731=\", \"77\u00d764=\", \"18\u00d773=\", \"31\u00d783=\", \"85\u00d747=\")\n  ,@(\"61\u00d737=\", \"44\u00d779=\", \"73\u00d755=\", \"59\u00d769=\", \"88\u00d767=\")\n  ,@(\"57\u00d787=\", \"27\u00d731=\", \"13\u00d799=\", \"68\u00d793=\", \"85\u00d766=\")\n  ,@(\"26\u00d730=\", \"46\u00d789=\", \"99\u00d757=\", \"62\u00d748=\", \"90\u00d782=\")\n  ,@(\"54\u00d791=\", \"77\u00d762=\", \"100\u00d761=\", \"50\u00d784=\", \"49\u00d742=\")\n  ,@(\"49\u00d714=\", \"93\u00d765=\", \"14\u00d718=\", \"22\u00d752=\", \"53\u00d757=\")\n  ,@(\"20\u00d748=\", \"92\u00d779=\", \"74\u00d722=\", \"44\u00d713=\", \"14\u00d759=\")\n  ,@(\"80\u00d790=\", \"73\u00d754=\", \"89\u00d786=\", \"20\u00d745=\", \"40\u00d736=\")\n  ,@(\"65\u00d777=\", \"45\u00d764=\", \"34\u00d756=\", \"33\u00d799=\", \"83\u00d759=\")\n  ,@(\"11\u00d744=\", \"35\u00d748=\", \"64\u00d732=\", \"93\u00d793=\", \"52\u00d723=\")\n)\n\nif ($t.Rows.Count -ne $newValues.Count -or $t.Columns.Count -ne $newValues[0].Count) {\n  throw \"Unexpected table shape: got $($t.Rows.Count)x$($t.Columns.Count), expected $($newValues.Count)x$($newValues[0].Count)\"\n}\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $t.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n  }\n}\n\n"}
